$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D-column cells whose new values look like plain numbers to remain text,
# matching the original inlineStr (text) cell type in the source workbook.
$textForceCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D16",
    "D17",
    "D19",
    "D21",
    "D22",
    "D23",
    "D24",
    "D26",
    "D27",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D40",
    "D41",
    "D43",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51",
)
foreach ($ref in $textForceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "41.201.94"
$ws.Range("E2").Value = "  -3.37%  "

$ws.Range("D3").Value = "2.460.76"
$ws.Range("E3").Value = "  -2.81%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "311.10"
$ws.Range("E5").Value = "  +0.48%  "

$ws.Range("D6").Value = "93.35"
$ws.Range("E6").Value = "  -6.69%  "

$ws.Range("D7").Value = "0.551"
$ws.Range("E7").Value = "  -2.90%  "

$ws.Range("D9").Value = "0.496"
$ws.Range("E9").Value = "  -4.92%  "

$ws.Range("D10").Value = "33.19"
$ws.Range("E10").Value = "  -7.12%  "

$ws.Range("D11").Value = "0.0776"
$ws.Range("E11").Value = "  -3.44%  "

$ws.Range("D12").Value = "0.108"
$ws.Range("E12").Value = "  -1.17%  "

$ws.Range("D13").Value = "6.96"
$ws.Range("E13").Value = "  -5.31%  "

$ws.Range("D14").Value = "2.839.08"

$ws.Range("D15").Value = "2.458.81"
$ws.Range("E15").Value = "  -4.51%  "

$ws.Range("D16").Value = "14.85"
$ws.Range("E16").Value = "  -3.03%  "

$ws.Range("D17").Value = "0.780"
$ws.Range("E17").Value = "  -3.99%  "

$ws.Range("D18").Value = "41.157.42"
$ws.Range("E18").Value = "  -3.46%  "

$ws.Range("D19").Value = "6.26"
$ws.Range("E19").Value = "  -6.85%  "

$ws.Range("D20").Value = "0.0₃0916"
$ws.Range("E20").Value = "  -3.52%  "

$ws.Range("D21").Value = "11.21"
$ws.Range("E21").Value = "  -8.76%  "

$ws.Range("D22").Value = "68.16"
$ws.Range("E22").Value = "  -1.60%  "

$ws.Range("D23").Value = "235.77"
$ws.Range("E23").Value = "  -3.06%  "

$ws.Range("D24").Value = "2.75"
$ws.Range("E24").Value = "  -4.01%  "

$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("D26").Value = "1.90"
$ws.Range("E26").Value = "  -6.22%  "

$ws.Range("D27").Value = "23.96"
$ws.Range("E27").Value = "  -6.24%  "

$ws.Range("E28").Value = "  -6.00%  "

$ws.Range("D29").Value = "9.57"
$ws.Range("E29").Value = "  -5.87%  "

$ws.Range("D30").Value = "35.94"
$ws.Range("E30").Value = "  -6.74%  "

$ws.Range("D31").Value = "152.43"
$ws.Range("E31").Value = "  -3.74%  "

$ws.Range("D32").Value = "5.46"
$ws.Range("E32").Value = "  -5.28%  "

$ws.Range("D33").Value = "2.65"
$ws.Range("E33").Value = "  -5.74%  "

$ws.Range("D34").Value = "2.53"
$ws.Range("E34").Value = "  -4.11%  "

$ws.Range("D35").Value = "0.0738"
$ws.Range("E35").Value = "  -5.59%  "

$ws.Range("D36").Value = "2.99"
$ws.Range("E36").Value = "  -4.61%  "

$ws.Range("D37").Value = "1.87"
$ws.Range("E37").Value = "  -4.95%  "

$ws.Range("D38").Value = "16.74"
$ws.Range("E38").Value = "  -8.88%  "

$ws.Range("E39").Value = "  -3.14%  "

$ws.Range("D40").Value = "0.101"
$ws.Range("E40").Value = "  -8.18%  "

$ws.Range("D41").Value = "4.17"
$ws.Range("E41").Value = "  -2.13%  "

$ws.Range("E42").Value = "  +0.20%  "

$ws.Range("D43").Value = "19.96"
$ws.Range("E43").Value = "  -11.25%  "

$ws.Range("D44").Value = "1.978.87"
$ws.Range("E44").Value = "  -0.39%  "

$ws.Range("D45").Value = "0.0283"
$ws.Range("E45").Value = "  -5.48%  "

$ws.Range("D46").Value = "3.02"
$ws.Range("E46").Value = "  -7.79%  "

$ws.Range("D47").Value = "8.68"
$ws.Range("E47").Value = "  -2.44%  "

$ws.Range("D48").Value = "69.03"
$ws.Range("E48").Value = "  -4.23%  "

$ws.Range("D49").Value = "96.63"
$ws.Range("E49").Value = "  -4.41%  "

$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "73.90"
$ws.Range("E50").Value = "  -7.03%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.176"
$ws.Range("E51").Value = "  -7.03%  "

